$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from header cell H1 (bold, centered, bordered) to new header cells I1 and J1,
# then set their text values.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "IF"

# Fill in the I and J data columns for rows 2-40.
$iValues = @(3,8,7,7,6,9,8,6,6,7,7,9,6,8,7,8,7,9,6,6,8,7,7,5,8,7,8,10,9,8,5,7,5,4,7,7,7,8,5)
$jValues = @(3,8,7,7,6,9,8,7,6,7,7,9,7,8,8,8,7,9,6,6,8,8,7,6,8,7,8,10,9,8,6,7,6,5,7,7,7,8,5)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}

$excel.CutCopyMode = $false
